$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '86.354.34'
$ws.Range("E2").Value = '  +8.50%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.314.74'
$ws.Range("E3").Value = '  +4.59%  '
$ws.Range("E4").Value = '  -0.18%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '217.50'
$ws.Range("E5").Value = '  +4.78%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '634.70'
$ws.Range("E6").Value = '  +1.17%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.319'
$ws.Range("E7").Value = '  +19.02%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.998'
$ws.Range("E8").Value = '  -0.24%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.594'
$ws.Range("E9").Value = '  -0.22%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '3.323.93'
$ws.Range("E10").Value = '  +4.54%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.598'
$ws.Range("E11").Value = '  -2.35%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.0000273'
$ws.Range("E12").Value = '  +6.92%  '
$ws.Range("E13").Value = '  +1.02%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '3.920.79'
$ws.Range("E14").Value = '  +3.72%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '33.99'
$ws.Range("E15").Value = '  +6.52%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '5.37'
$ws.Range("E16").Value = '  +1.59%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '85.918.83'
$ws.Range("E17").Value = '  +7.43%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '3.312.08'
$ws.Range("E18").Value = '  +3.62%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '14.56'
$ws.Range("E19").Value = '  +1.32%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '3.15'
$ws.Range("E20").Value = '  +7.23%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '442.79'
$ws.Range("E21").Value = '  +0.48%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '9.11'
$ws.Range("E22").Value = '  -0.60%  '
$ws.Range("E23").Value = '  -1.56%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '7.42'
$ws.Range("E24").Value = '  +5.83%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '5.39'
$ws.Range("E25").Value = '  +14.54%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '12.18'
$ws.Range("E26").Value = '  +12.08%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '3.491.03'
$ws.Range("E27").Value = '  +3.24%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '78.12'
$ws.Range("E28").Value = '  +1.64%  '
$ws.Range("E29").Value = '  +5.70%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.998'
$ws.Range("E30").Value = '  -0.24%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.169'
$ws.Range("E31").Value = '  +38.35%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '598.80'
$ws.Range("E32").Value = '  +9.72%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '9.18'
$ws.Range("E33").Value = '  +0.93%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.999'
$ws.Range("E34").Value = '  -0.24%  '
$ws.Range("E35").Value = '  +4.12%  '
$ws.Range("E36").Value = '  +1.78%  '
$ws.Range("E37").Value = '  +0.59%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '23.31'
$ws.Range("E38").Value = '  +0.33%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '6.45'
$ws.Range("E39").Value = '  +14.06%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.998'
$ws.Range("E40").Value = '  +0.08%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.416'
$ws.Range("E41").Value = '  +1.53%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '21.31'
$ws.Range("E42").Value = '  +2.60%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '3.11'
$ws.Range("E43").Value = '  +15.70%  '
$ws.Range("E44").Value = '  +13.02%  '
$ws.Range("E45").Value = '  -0.01%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '157.78'
$ws.Range("E46").Value = '  -4.17%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '187.41'
$ws.Range("E47").Value = '  -0.90%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.36'
$ws.Range("E48").Value = '  +4.02%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '45.27'
$ws.Range("E49").Value = '  +3.75%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.788'
$ws.Range("E50").Value = '  +0.63%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '26.17'
$ws.Range("E51").Value = '  +4.89%  '

Write-Host "Applied 91 cell updates"
